$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "66.746.56"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "3.488.18"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.67%  "
$ws.Range("D7").Value = "3.486.26"
$ws.Range("E7").Value = "  -0.98%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.483"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.38%  "
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000215"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.53%  "
$ws.Range("D14").Value = "4.073.63"
$ws.Range("E14").Value = "  -1.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.25%  "
$ws.Range("D16").Value = "3.478.52"
$ws.Range("E16").Value = "  -2.15%  "
$ws.Range("D17").Value = "66.826.30"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.87%  "
$ws.Range("E22").Value = "  -2.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.612"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.34%  "
$ws.Range("D26").Value = "3.623.57"
$ws.Range("E26").Value = "  -1.77%  "
$ws.Range("E27").Value = "  -6.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.80"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.51%  "
$ws.Range("E29").Value = "  -5.72%  "
$ws.Range("E30").Value = "  -2.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.170"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.45"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.36%  "
$ws.Range("D36").Value = "3.476.02"
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("E37").Value = "  -4.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.94"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.74%  "
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "176.82"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.894"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.06"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("E48").Value = "  -5.88%  "
$ws.Range("E49").Value = "  -3.06%  "
$ws.Range("E50").Value = "  -7.55%  "
$ws.Range("E51").Value = "  -2.45%  "
